$wb = $excel.ActiveWorkbook

$struct = $wb.Worksheets.Item("strucutres")

# Insert a new first column on the "strucutres" sheet, shifting all
# existing data (and their styles) one column to the right.
$struct.Columns.Item(1).Insert()

# Populate the new "Status" column with its header and the two notes.
$struct.Range("A1").Value = "Status"
$struct.Range("A2").Value = "made synonym request on phenotype-ext"
$struct.Range("A3").Value = "on fovt-humerus; made pull request"

# Move the active/selected tab from "axis" to "strucutres", and update
# the selected cell on "strucutres" to A3.
$struct.Activate()
$struct.Range("A3").Select()
